$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two new numbered columns (10, 11) ---
# Reuse the existing header cell style (bold, centered, bordered) by copying
# format from K1 (which already carries that style) onto the new L1:M1 cells.
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11

# --- Add a new row 7, matching the same style used by the other index cells in column A ---
$ws.Range("A5").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# --- Row index column (A) values ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

# --- Force Text format on the data columns so numeric-looking strings (section
#     numbers, CRNs, times, leading zeros, etc.) are stored as text, not numbers ---
$ws.Range("B2:M7").NumberFormat = "@"

$ws.Range("B2").Value = "102"
$ws.Range("C2").Value = "01"
$ws.Range("D2").Value = "7273"
$ws.Range("E2").Value = "4"
$ws.Range("F2").Value = "F"
$ws.Range("G2").Value = "1000"
$ws.Range("H2").Value = "Beginning Italian II MTWRF"
$ws.Range("I2").Value = "1050am"
$ws.Range("J2").Value = "PETE"
$ws.Range("K2").Value = "302"
$ws.Range("L2").Value = "Di"
$ws.Range("M2").Value = "Siena Ivana"
$ws.Range("B3").Value = "102"
$ws.Range("C3").Value = "02"
$ws.Range("D3").Value = "11406"
$ws.Range("E3").Value = "4"
$ws.Range("F3").Value = "F"
$ws.Range("G3").Value = "1100"
$ws.Range("H3").Value = "Beginning Italian II MTWRF"
$ws.Range("I3").Value = "1150am"
$ws.Range("J3").Value = "PETE"
$ws.Range("K3").Value = "302"
$ws.Range("L3").Value = "Di"
$ws.Range("M3").Value = "Siena Ivana"
$ws.Range("B4").Value = "301"
$ws.Range("C4").Value = "01"
$ws.Range("D4").Value = "19401"
$ws.Range("E4").Value = "2"
$ws.Range("F4").Value = "F"
$ws.Range("G4").Value = "MW"
$ws.Range("H4").Value = "Intro to Italian Literature"
$ws.Range("I4").Value = "0130"
$ws.Range("J4").Value = "0220pm"
$ws.Range("K4").Value = "PETE"
$ws.Range("L4").Value = "102"
$ws.Range("M4").Value = "Milkova Stiliana"
$ws.Range("B5").Value = "401"
$ws.Range("C5").Value = "01"
$ws.Range("D5").Value = "19402"
$ws.Range("E5").Value = "2"
$ws.Range("F5").Value = "F"
$ws.Range("G5").Value = "MW"
$ws.Range("H5").Value = "Contemporary Italian Lit"
$ws.Range("I5").Value = "0130"
$ws.Range("J5").Value = "0220pm"
$ws.Range("K5").Value = "PETE"
$ws.Range("L5").Value = "102"
$ws.Range("M5").Value = "Milkova Stiliana"
$ws.Range("B6").Value = "995F"
$ws.Range("C6").Value = "01"
$ws.Range("D6").Value = "16969"
$ws.Range("E6").Value = "4"
$ws.Range("F6").Value = "F"
$ws.Range("H6").Value = "Private Reading "
$ws.Range("I6").Value = "Full"
$ws.Range("J6").Value = "TBA"
$ws.Range("K6").Value = "TBA"
$ws.Range("L6").Value = "Di"
$ws.Range("M6").Value = "Siena Ivana"
$ws.Range("B7").Value = "995H"
$ws.Range("C7").Value = "01"
$ws.Range("D7").Value = "16968"
$ws.Range("E7").Value = "2"
$ws.Range("F7").Value = "F"
$ws.Range("H7").Value = "Private Reading "
$ws.Range("I7").Value = "Half"
$ws.Range("J7").Value = "TBA"
$ws.Range("K7").Value = "TBA"
$ws.Range("L7").Value = "Di"
$ws.Range("M7").Value = "Siena Ivana"

# --- Restore the plain/default formatting (no explicit style) on the data
#     cells, matching the rest of the table, now that the text values are set ---
$ws.Range("A1").Copy()
$ws.Range("B2:M7").PasteSpecial(-4122)

# --- Drop the day-of-week cell for the two TBA "Private Reading" rows ---
$ws.Range("G6").Value = $null
$ws.Range("G7").Value = $null
